# Update Test_Results sheet with latest automated test run outcomes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Step 1 (Open URL) now passes
$ws.Range("L2").Value = "PASS"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = "Opened: https://10.0.49.147:7504/FCJNeoWeb/SMMDIFRM.jsp"
$ws.Range("O2").Value = "screenshots/STEP_1.png"
$ws.Range("P2").Value = "page_sources/STEP_1_source.html"

# Rows that were skipped because "TO BE EXECUTED" = NO for their scenario
$skippedRows = @(3, 4, 5, 6, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40)
foreach ($r in $skippedRows) {
    $ws.Range("L$r").Value = "SKIPPED"
    $ws.Range("M$r").Value = "TO BE EXECUTED = NO"
}

# Row 7: Step 6 (Fill Function Id) passes
$ws.Range("L7").Value = "PASS"
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = "Filled: Function Id"
$ws.Range("O7").Value = "screenshots/STEP_6.png"
$ws.Range("P7").Value = "page_sources/STEP_6_source.html"

# Row 8: Step 7 (Click Go) passes
$ws.Range("L8").Value = "PASS"
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = "Clicked: Go"
$ws.Range("O8").Value = "screenshots/STEP_7.png"
$ws.Range("P8").Value = "page_sources/STEP_7_source.html"

# Row 9: Step 8 (Click New) fails - could not click element
$ws.Range("L9").Value = "FAIL"
$ws.Range("M9").Value = "Could not click element"
$ws.Range("N9").Value = "Failed to click: New"
$ws.Range("O9").Value = ""
$ws.Range("P9").Value = ""

# Rows 10-24: subsequent steps fail because no valid page was available
$noPageRows = @(10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24)
foreach ($r in $noPageRows) {
    $ws.Range("L$r").Value = "FAIL"
    $ws.Range("M$r").Value = "No valid page available"
    $ws.Range("N$r").Value = "No valid page available"
    $ws.Range("O$r").Value = ""
    $ws.Range("P$r").Value = ""
}
